# Apply the "Automatic update of files" edit to the Artfynd sheet.
# Rows 48-53 get their species-observation data re-shuffled/updated as per
# the upstream source refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# --- Row 48 ---
$ws.Range("A48").Value = 111974134
$ws.Range("B48").Value = 90792
$ws.Range("E48").Value = 4361
$ws.Range("F48").Value = "Orange taggsvamp"
$ws.Range("G48").Value = "Hydnellum aurantiacum"
$ws.Range("H48").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("Q48").Value = 439400
$ws.Range("R48").Value = 6952207

# --- Row 49 ---
$ws.Range("A49").Value = 111974029
$ws.Range("B49").Value = 88166
$ws.Range("D49").Value = "VU"
$ws.Range("E49").Value = 6276
$ws.Range("F49").Value = "Goliatmusseron"
$ws.Range("G49").Value = "Tricholoma matsutake"
$ws.Range("H49").Value = "(S.Ito & S.Imai) Singer"
$ws.Range("P49").Value = "Aloppmoarna, Jmt"
$ws.Range("Q49").Value = 439335
$ws.Range("R49").Value = 6952297

# --- Row 50 ---
$ws.Range("A50").Value = 111974124
$ws.Range("B50").Value = 90800
$ws.Range("D50").Value = "LC"
$ws.Range("E50").Value = 4364
$ws.Range("F50").Value = "Dropptaggsvamp"
$ws.Range("G50").Value = "Hydnellum ferrugineum"
$ws.Range("H50").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q50").Value = 439276
$ws.Range("R50").Value = 6952197

# --- Row 51 ---
$ws.Range("A51").Value = 111974126
$ws.Range("B51").Value = 88166
$ws.Range("P51").Value = "Aloppmoarna i S, Jmt"
$ws.Range("Q51").Value = 439290
$ws.Range("R51").Value = 6952209

# --- Row 52 ---
$ws.Range("B52").Value = 90794

# --- Row 53 ---
$ws.Range("A53").Value = 111974133
$ws.Range("B53").Value = 90816
$ws.Range("D53").Value = "NT"
$ws.Range("E53").Value = 2059
$ws.Range("F53").Value = "Skrovlig taggsvamp"
$ws.Range("G53").Value = "Hydnellum scabrosum"
$ws.Range("H53").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("Q53").Value = 439390
$ws.Range("R53").Value = 6952220
